# Updates the "Price" (D) and "Volume(1h)" (E) columns of the cryptos
# sheet with freshly scraped values.
#
# Note: a handful of new Price values (e.g. "23.44") look like plain
# decimal numbers. A bare `.Value = "23.44"` assignment would be
# auto-converted by Excel into a numeric cell (23.44 as a float),
# which does not match the source data (these columns are text).
# To force those specific cells to stay text -- exactly as typing
# '23.44 into Excel does -- a leading apostrophe (quote-prefix) is
# used for those values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.913.19'
$ws.Range('E2').Value = '  -0.07%  '
$ws.Range('D3').Value = '1.635.56'
$ws.Range('E3').Value = '  +0.16%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = '''211.75'
$ws.Range('E5').Value = '  -0.09%  '
$ws.Range('E6').Value = '  -0.61%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('D8').Value = '''23.44'
$ws.Range('E8').Value = '  +0.68%  '
$ws.Range('E9').Value = '  -0.39%  '
$ws.Range('E10').Value = '  -0.41%  '
$ws.Range('E11').Value = '  +0.17%  '
$ws.Range('D12').Value = '1.869.34'
$ws.Range('E12').Value = '  +0.27%  '
$ws.Range('D13').Value = '1.643.14'
$ws.Range('E13').Value = '  +0.58%  '
$ws.Range('D14').Value = '''0.564'
$ws.Range('E14').Value = '  -0.26%  '
$ws.Range('E15').Value = '  -0.98%  '
$ws.Range('D16').Value = '''65.33'
$ws.Range('E16').Value = '  +0.07%  '
$ws.Range('D17').Value = '27.941.90'
$ws.Range('D18').Value = '''230.35'
$ws.Range('E18').Value = '  -0.17%  '
$ws.Range('D19').Value = '''7.84'
$ws.Range('E19').Value = '  +4.10%  '
$ws.Range('D20').Value = '0.0₃0719'
$ws.Range('E20').Value = '  -0.37%  '
$ws.Range('E21').Value = '  -0.08%  '
$ws.Range('E22').Value = '  -0.02%  '
$ws.Range('D23').Value = '''10.17'
$ws.Range('E23').Value = '  -2.11%  '
$ws.Range('E24').Value = '  +0.99%  '
$ws.Range('D25').Value = '''156.44'
$ws.Range('E25').Value = '  +1.13%  '
$ws.Range('E26').Value = '  +0.06%  '
$ws.Range('E27').Value = '  -0.02%  '
$ws.Range('E28').Value = '  -0.48%  '
$ws.Range('E29').Value = '  +0.13%  '
$ws.Range('D30').Value = '''1.18'
$ws.Range('E30').Value = '  +0.00%  '
$ws.Range('E31').Value = '  -0.21%  '
$ws.Range('E32').Value = '  +0.98%  '
$ws.Range('D33').Value = '''3.11'
$ws.Range('E33').Value = '  +1.03%  '
$ws.Range('D34').Value = '1.404.63'
$ws.Range('E34').Value = '  +0.11%  '
$ws.Range('E35').Value = '  +3.07%  '
$ws.Range('E36').Value = '  +0.80%  '
$ws.Range('E37').Value = '  -0.80%  '
$ws.Range('E38').Value = '  +0.60%  '
$ws.Range('D39').Value = '''0.560'
$ws.Range('E39').Value = '  -0.39%  '
$ws.Range('D40').Value = '''0.853'
$ws.Range('E40').Value = '  -2.21%  '
$ws.Range('E41').Value = '  -0.03%  '
$ws.Range('D42').Value = '''1.00'
$ws.Range('E42').Value = '  -1.97%  '
$ws.Range('E43').Value = '  +1.91%  '
$ws.Range('D44').Value = '''66.15'
$ws.Range('E44').Value = '  -1.15%  '
$ws.Range('D45').Value = '''5.46'
$ws.Range('E45').Value = '  -1.50%  '
$ws.Range('D46').Value = '1.777.66'
$ws.Range('E46').Value = '  +0.16%  '
$ws.Range('D47').Value = '''2.14'
$ws.Range('E47').Value = '  -2.45%  '
$ws.Range('D48').Value = '''88.65'
$ws.Range('E48').Value = '  +1.00%  '
$ws.Range('E49').Value = '  +2.45%  '
$ws.Range('E50').Value = '  -0.20%  '
$ws.Range('D51').Value = '''7.68'
$ws.Range('E51').Value = '  +2.94%  '
